# Update crypto price/volume data per diff (Fri Mar 22 21:12:40 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are stored as text, matching source data
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '63.222.29'
$ws.Cells.Item(2, 5).Value = '  -3.38%  '
$ws.Cells.Item(3, 4).Value = '3.303.23'
$ws.Cells.Item(3, 5).Value = '  -5.09%  '
$ws.Cells.Item(4, 5).Value = '  +0.13%  '
$ws.Cells.Item(5, 4).Value = '544.98'
$ws.Cells.Item(5, 5).Value = '  -1.46%  '
$ws.Cells.Item(6, 4).Value = '170.25'
$ws.Cells.Item(6, 5).Value = '  -4.67%  '
$ws.Cells.Item(7, 4).Value = '0.612'
$ws.Cells.Item(7, 5).Value = '  -4.05%  '
$ws.Cells.Item(8, 5).Value = '  +0.02%  '
$ws.Cells.Item(9, 2).Value = 'LidoStakedEther'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Cells.Item(9, 4).Value = '3.291.05'
$ws.Cells.Item(9, 5).Value = '  -5.29%  '
$ws.Cells.Item(10, 2).Value = 'Cardano'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(10, 4).Value = '0.606'
$ws.Cells.Item(10, 5).Value = '  -4.08%  '
$ws.Cells.Item(11, 2).Value = 'Dogecoin'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(11, 4).Value = '0.150'
$ws.Cells.Item(11, 5).Value = '  -0.58%  '
$ws.Cells.Item(12, 2).Value = 'Avalanche'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(12, 4).Value = '52.69'
$ws.Cells.Item(12, 5).Value = '  -1.76%  '
$ws.Cells.Item(13, 2).Value = 'ShibaInu'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(13, 4).Value = '0.0000263'
$ws.Cells.Item(13, 5).Value = '  -2.58%  '
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).Value = '8.79'
$ws.Cells.Item(14, 5).Value = '  -4.66%  '
$ws.Cells.Item(15, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(15, 4).Value = '3.844.28'
$ws.Cells.Item(15, 5).Value = '  -4.82%  '
$ws.Cells.Item(16, 2).Value = 'Chainlink'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(16, 4).Value = '17.92'
$ws.Cells.Item(16, 5).Value = '  -3.93%  '
$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = '3.312.28'
$ws.Cells.Item(17, 5).Value = '  -4.83%  '
$ws.Cells.Item(18, 2).Value = 'TRON'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(18, 4).Value = '0.116'
$ws.Cells.Item(18, 5).Value = '  -4.12%  '
$ws.Cells.Item(19, 2).Value = 'WrappedBTC'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(19, 4).Value = '63.400.54'
$ws.Cells.Item(19, 5).Value = '  -3.17%  '
$ws.Cells.Item(20, 2).Value = 'Uniswap'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(20, 4).Value = '11.57'
$ws.Cells.Item(20, 5).Value = '  -3.74%  '
$ws.Cells.Item(21, 2).Value = 'Polygon'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(21, 4).Value = '0.966'
$ws.Cells.Item(21, 5).Value = '  -2.13%  '
$ws.Cells.Item(22, 2).Value = 'BitcoinCash'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(22, 4).Value = '411.02'
$ws.Cells.Item(22, 5).Value = '  -1.21%  '
$ws.Cells.Item(23, 2).Value = 'PancakeSwap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(23, 4).Value = '4.01'
$ws.Cells.Item(23, 5).Value = '  -0.68%  '
$ws.Cells.Item(24, 2).Value = 'Toncoin'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(24, 4).Value = '4.36'
$ws.Cells.Item(24, 5).Value = '  +6.03%  '
$ws.Cells.Item(25, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(25, 4).Value = '13.62'
$ws.Cells.Item(25, 5).Value = '  +5.05%  '
$ws.Cells.Item(26, 2).Value = 'Litecoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(26, 4).Value = '82.46'
$ws.Cells.Item(26, 5).Value = '  -4.26%  '
$ws.Cells.Item(27, 2).Value = 'RenderToken'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(27, 4).Value = '10.45'
$ws.Cells.Item(27, 5).Value = '  -2.97%  '
$ws.Cells.Item(28, 2).Value = 'ImmutableX'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(28, 4).Value = '2.70'
$ws.Cells.Item(28, 5).Value = '  -5.13%  '
$ws.Cells.Item(29, 2).Value = 'Filecoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(29, 4).Value = '8.54'
$ws.Cells.Item(29, 5).Value = '  -6.02%  '
$ws.Cells.Item(30, 2).Value = 'EthereumClassic'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(30, 4).Value = '28.80'
$ws.Cells.Item(30, 5).Value = '  -4.70%  '
$ws.Cells.Item(31, 2).Value = 'NEARProtocol'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(31, 4).Value = '6.30'
$ws.Cells.Item(31, 5).Value = '  -3.32%  '
$ws.Cells.Item(32, 2).Value = 'Cosmos'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(32, 4).Value = '11.24'
$ws.Cells.Item(32, 5).Value = '  -4.21%  '
$ws.Cells.Item(33, 2).Value = 'Bittensor'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(33, 4).Value = '567.47'
$ws.Cells.Item(33, 5).Value = '  -7.38%  '
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 4).Value = '0.105'
$ws.Cells.Item(34, 5).Value = '  -3.79%  '
$ws.Cells.Item(35, 2).Value = 'OKB'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(35, 4).Value = '57.50'
$ws.Cells.Item(35, 5).Value = '  -2.85%  '
$ws.Cells.Item(36, 2).Value = 'Dai'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(36, 4).Value = '1.00'
$ws.Cells.Item(36, 5).Value = '  -0.12%  '
$ws.Cells.Item(37, 2).Value = 'Kaspa'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(37, 4).Value = '0.145'
$ws.Cells.Item(37, 5).Value = '  -0.48%  '
$ws.Cells.Item(38, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(38, 4).Value = '34.79'
$ws.Cells.Item(38, 5).Value = '  -6.83%  '
$ws.Cells.Item(39, 2).Value = 'Stacks'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(39, 4).Value = '3.37'
$ws.Cells.Item(39, 5).Value = '  +4.16%  '
$ws.Cells.Item(40, 2).Value = 'PEPE'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(40, 4).Value = '0.0₃0728'
$ws.Cells.Item(40, 5).Value = '  -6.87%  '
$ws.Cells.Item(41, 2).Value = 'TheGraph'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(41, 4).Value = '0.362'
$ws.Cells.Item(41, 5).Value = '  -4.46%  '
$ws.Cells.Item(42, 2).Value = 'Maker'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(42, 4).Value = '3.104.65'
$ws.Cells.Item(42, 5).Value = '  -8.05%  '
$ws.Cells.Item(43, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(43, 4).Value = '1.00'
$ws.Cells.Item(43, 5).Value = '  +0.19%  '
$ws.Cells.Item(44, 2).Value = 'ThetaToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(44, 4).Value = '2.74'
$ws.Cells.Item(44, 5).Value = '  -2.25%  '
$ws.Cells.Item(45, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(45, 4).Value = '3.19'
$ws.Cells.Item(45, 5).Value = '  -1.29%  '
$ws.Cells.Item(46, 2).Value = 'VeChain'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(46, 4).Value = '0.0397'
$ws.Cells.Item(46, 5).Value = '  -4.08%  '
$ws.Cells.Item(47, 2).Value = 'Fetch.AI'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(47, 4).Value = '2.39'
$ws.Cells.Item(47, 5).Value = '  -5.42%  '
$ws.Cells.Item(48, 2).Value = 'WEMIXToken'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(48, 4).Value = '2.59'
$ws.Cells.Item(48, 5).Value = '  -4.33%  '
$ws.Cells.Item(49, 2).Value = 'Stellar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(49, 4).Value = '0.127'
$ws.Cells.Item(49, 5).Value = '  -3.84%  '
$ws.Cells.Item(50, 2).Value = 'Monero'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(50, 4).Value = '133.68'
$ws.Cells.Item(50, 5).Value = '  -2.92%  '
$ws.Cells.Item(51, 2).Value = 'THORChain'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(51, 4).Value = '7.96'
$ws.Cells.Item(51, 5).Value = '  -5.66%  '
